# Generate Report for Handoff
#
# The localization-status report is regenerated for a new handoff:
#  - the "Status" column (previously "Handed back: in sync with en-US")
#    becomes "Ready for handoff" on every sheet (Overview, zh-cn, de-de);
#  - the associated handoff/generation timestamps are refreshed;
#  - because "Ready for handoff" is shorter than the old status text, the
#    status/date columns are narrowed to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status + HO Xliff generation date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-22 03:08:41"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-22 03:08:37"

# --- de-de detail sheet: Status + Latest Handoff Datetime ---------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-22 03:08:41"

# --- Narrow the status/date columns now that the text is shorter -------
# (target character width ~17.216; the host's ColumnWidth setter quantizes
# to the nearest 1/6 of a character, so 98/6 = 16.333... in is the closest
# input that lands on the nearest achievable width.)
$overview.Range("E1").ColumnWidth = 16.3333333333333
$overview.Range("F1").ColumnWidth = 16.3333333333333
$zhcn.Range("C1").ColumnWidth = 16.3333333333333
$dede.Range("C1").ColumnWidth = 16.3333333333333
